# Nexial "number-showcase" workbook update
#
# [web] - deselect(locator,text): NEW command to deselect from a SELECT
#         element based on `text`, inserted into the hidden '#system' sheet's
#         "web" reference list (column U) at its correct alphabetical
#         position (right before "deselectMulti(locator,array)").
#
# This requires:
#   1. Shifting the existing column-U values (rows 53..116) down by one row
#      (to rows 54..117), without disturbing any of the other independent
#      reference-list columns (A..Z) that happen to share the same row
#      numbers on this sheet.
#   2. Writing the new command text into the now-empty U53.
#   3. Growing the "web" defined name from $U$2:$U$116 to $U$2:$U$117 so the
#      dropdown/data-validation list picks up the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$firstRow = 53
$lastRowBefore = 116
$lastRowAfter = $lastRowBefore + 1

# Shift column U values down by one row, starting from the bottom so we
# never overwrite a value before it has been copied down.
for ($r = $lastRowBefore; $r -ge $firstRow; $r--) {
    $srcCell = $ws.Range("U" + $r)
    $dstCell = $ws.Range("U" + ($r + 1))
    $dstCell.Value = $srcCell.Value2
}

# Insert the new command text in its correct alphabetical position.
$ws.Range("U" + $firstRow).Value = "deselect(locator,text)"

# Update the "web" defined name so it covers the new row.
$webName = $wb.Names.Item("web")
$webName.RefersTo = "='#system'!`$U`$2:`$U$" + $lastRowAfter
